$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (row 7) - columns M and N
$ws.Range("M7").Value = "AlarmLoadingDetail"
$ws.Range("N7").Value = "StandbyLoadingDetail"

# New data values (rows 8 and 9) - columns M and N
$ws.Range("M8").Value = "Battery Alarm (A)"
$ws.Range("N8").Value = "Battery Standby (A)"

$ws.Range("M9").Value = "Battery Alarm (A)"
$ws.Range("N9").Value = "Battery Standby (A)"

# Match formats/styles used on similar cells (reuse existing style entries
# rather than creating new ones) by copying formats only.
$ws.Range("A7").Copy()
$ws.Range("M7:N7").PasteSpecial(-4122)

$ws.Range("I8").Copy()
$ws.Range("M8:N8").PasteSpecial(-4122)

$ws.Range("I9").Copy()
$ws.Range("M9:N9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update selection to match diff
$ws.Range("M9:N9").Select()
